$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new work-log entry as row 31
$ws.Range("A31").Value2 = 43490
$ws.Range("A31").NumberFormat = "m/d/yy"
$ws.Range("B31").Value = "Updated the Methodologies section "
$ws.Range("C31").Value = 1.75

# Update the view to match where Excel scrolled/selected after the edit
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("E33").Select()
